# Fix spreadsheet parsing errors
# 1) Remove the "Diana" matchup row from the "Mid" sheet (row 12), shifting
#    all rows below it up by one.
# 2) Change Ahri's difficulty value from the text "1 or 2" to the numeric 1.5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mid")

# Find and delete the row whose Champion column (A) is "Diana".
$found = $ws.Cells.Find("Diana")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}

# Ahri's difficulty cell is column B on the same row as "Ahri" in column A.
$ahriCell = $ws.Cells.Find("Ahri")
$diffCell = $ws.Cells.Item($ahriCell.Row, 2)
$diffCell.Value = 1.5

# Make the "Mid" sheet the active sheet / tab, matching the saved selection.
$ws.Activate()
$ws.Range("C4").Select()

$meta = $wb.Worksheets.Item("metadata")
$meta.Range("B7").Select()
